$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "end_year"/"start_year" headers to "end"/"start"
$ws.Range("D1").Value = "end"
$ws.Range("K1").Value = "start"

# Move the active selection from A3 to L3
$ws.Range("L3").Select() | Out-Null
